$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("soknad")

# Three new kindergarten applications (soknad) were submitted. Insert 3 blank
# rows right after the header row, pushing the existing applications down,
# then renumber the id (A, zero-based row counter) and sok_id (B, countdown
# counter) columns for the whole data range so they stay consistent, and
# restore the inserted rows' formatting to match the rest of the table.
$ws.Rows("2:4").Insert()

# Column A keeps the bordered/bold "index" style used throughout the table;
# the rest of the new rows should have no special formatting, same as every
# other data row.
$ws.Range("A5").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)
$ws.Range("B2:M4").ClearFormats()
$excel.CutCopyMode = $false

$lastRow = 14
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
    $ws.Cells.Item($r, 2).Value = $lastRow + 1 - $r
}
